$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted values in the source
# workbook (e.g. "53.90", "0.140", thousand-dot separators, percent strings).
# Force each cell to Text format immediately before writing it so Excel does not
# reinterpret/round the string as a number (applying NumberFormat via one big
# multi-area Range only reliably sticks to the first area, so it is done here
# per-cell right before the corresponding value write).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.305.38'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.846.75'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '360.89'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '113.73'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.576'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.609'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.66'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0865'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.08'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.132'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.83'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.291.43'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.848.81'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.908'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '52.197.39'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.62'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +9.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.16'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.59'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0996'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.43'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.74'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.23'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.44'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '53.90'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.141'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'VeChain'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0465'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +23.26%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.71'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.90'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.42'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +8.09%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.08'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.47'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.96'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.117'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '128.11'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.02%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.17%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.43'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.122.94'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.96%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +10.33%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.50'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.94%  '
